# Update countries & provincias Spain
# - Refresh "datos actualizados" timestamp
# - Martinica's case count overtook Lesoto's -> rows swap order (sorted desc by Casos totales)
# - Nueva Caledonia / Santa Lucia tied -> swap order (cosmetic, same values)
# - Refresh case numbers for several countries (new scrape)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Timestamp banner in A1
$ws.Range("A1").Value = "Datos actualizados a 7 de Octubre de 2020 a las 03:23"

# 2) Estados Unidos (row 4) - refreshed totals
$ws.Range("B4").Value = 7719403
$ws.Range("C4").Value = 40317
$ws.Range("D4").Value = 4934770
$ws.Range("E4").Value = 2568853
$ws.Range("G4").Value = 748
$ws.Range("H4").Value = 215780

# 3) Alemania (row 26)
$ws.Range("D26").Value = 265600
$ws.Range("E26").Value = 31884

# 4) Belgica (row 35)
$ws.Range("D35").Value = 19719
$ws.Range("E35").Value = 102406

# 5) Panama (row 37)
$ws.Range("B37").Value = 116602
$ws.Range("C37").Value = 683
$ws.Range("D37").Value = 92950
$ws.Range("E37").Value = 21212
$ws.Range("G37").Value = 10
$ws.Range("H37").Value = 2440

# 6) Venezuela (row 56)
$ws.Range("B56").Value = 79796
$ws.Range("C56").Value = 679
$ws.Range("D56").Value = 70719
$ws.Range("E56").Value = 8412
$ws.Range("G56").Value = 7
$ws.Range("H56").Value = 665

# 7) Australia (row 81)
$ws.Range("B81").Value = 27174
$ws.Range("C81").Value = 25
$ws.Range("D81").Value = 24917
$ws.Range("E81").Value = 1362

# 8) Lesoto / Martinica reorder (row 163 was Lesoto, row 164 was Martinica).
#    Martinica's new totals push it above Lesoto, so row 163 becomes Martinica
#    (with the fresh numbers) and row 164 becomes Lesoto (keeping its prior numbers).
$ws.Range("A163").Value = "Martinica"
$ws.Range("B163").Value = 1851
$ws.Range("C163").Value = 308
$ws.Range("D163").Value = 98
$ws.Range("E163").Value = 1731
$ws.Range("F163").Value = 0
$ws.Range("G163").Value = 1
$ws.Range("H163").Value = 22

$ws.Range("A164").Value = "Lesoto"
$ws.Range("B164").Value = 1683
$ws.Range("C164").Value = 0
$ws.Range("D164").Value = 926
$ws.Range("E164").Value = 718
$ws.Range("F164").Value = 0
$ws.Range("G164").Value = 0
$ws.Range("H164").Value = 39

# 9) Santo Tome y Principe (row 169)
$ws.Range("B169").Value = 914
$ws.Range("C169").Value = 1
$ws.Range("E169").Value = 11

# 10) Mauricio (row 184)
$ws.Range("B184").Value = 395
$ws.Range("C184").Value = 8
$ws.Range("E184").Value = 28

# 11) Santa Lucia / Nueva Caledonia reorder (row 207 was Santa Lucia, row 208 was
#     Nueva Caledonia). Same totals for both, just swap the country labels so
#     Nueva Caledonia sorts first.
$ws.Range("A207").Value = "Nueva Caledonia"
$ws.Range("A208").Value = "Santa Lucia"
